$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GL")
$ws.Activate()

# --- Insert a fresh row at 46 (shifts the former row 46 - date 42783 / TV note - down to 47).
# Excel gives the new blank row the "format same as the row above" styles (date style
# from A45, #,##0 style from B45) for free, which is exactly what row 46 needs. ---
$ws.Rows.Item(46).Insert()

# --- Insert 8 more rows below the (now relocated) row 47, again picking up "format same
# as above" -- this time the d-mmm date style that A47 carries -- for rows 48-55. ---
$ws.Rows.Item(48).Resize(8).Insert()

# Row 1's two running totals get stretched by these inserts; pin them back to the exact
# one-row expansion the workbook actually ended up with (SUM(C2:C197)->SUM(C2:C198),
# SUM(B2:B198)->SUM(B2:B199)) rather than the ten rows' worth Excel's auto-adjust applied.
$ws.Range("E1").Formula = "=SUM(C2:C198)"
$ws.Range("H1").Formula = "=SUM(B2:B199)"

# New row 46: a "Quincena" entry dated 2/14/2017.
$ws.Range("A46").Value = 42780
$ws.Range("B46").Value = 2500
$ws.Range("C46").Value = 125.35
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "Quincena"

# Row 47 is the old row 46 (A47/G47 already correct: 42783 / "Mande $6 mil..."). Add the
# Television / regalo details that go with that TV purchase.
$ws.Range("B47").Value = 6400
$ws.Range("C47").Value = 482.45
$ws.Range("D47").Value = "Television"
$ws.Range("E47").Formula = "=B47/C47"

$ws.Range("A48").Value = 42783
$ws.Range("B48").Value = -4000
$ws.Range("C48").Formula = "=B48/E47"
$ws.Range("D48").Value = "Regalo Lili"
$ws.Range("G48").Clear()

$ws.Range("A49").Value = 42793
$ws.Range("B49").Value = 2500
$ws.Range("C49").Value = 129
$ws.Range("D49").Value = "Quincena"
$ws.Range("G49").Clear()

$ws.Range("A50").Value = 42804
$ws.Range("B50").Value = 1600
$ws.Range("C50").Value = 82.75
$ws.Range("D50").Value = "Doctor visit and medicines"
$ws.Range("G50").Clear()

$ws.Range("A51").Value = 42809
$ws.Range("B51").Value = 2500
$ws.Range("C51").Value = 130
$ws.Range("D51").Value = "Quincena"
$ws.Range("G51").Clear()

$ws.Range("A52").Value = 42824
$ws.Range("B52").Value = 4100
$ws.Range("C52").Value = 223
$ws.Range("D52").Value = "Quincena + Dr visit and medicines ($1,600)"
$ws.Range("G52").Clear()

# Row 56 only carries the running note in column G -- written here (ahead of the two
# "Renta casa Lili"/"Medicinas Laila" rows below) so the shared-string table fills up
# in the same order the original author typed things in.
$ws.Range("G56").Value = "Mande $4,500; $2,500 quincena, + $2,000 cambio de casa"
$ws.Range("G56").Interior.Color = $ws.Range("G47").Interior.Color

$ws.Range("A53").Value = 42828
$ws.Range("B53").Value = 9000
$ws.Range("C53").Value = 489
$ws.Range("D53").Value = "Renta casa Lili, + seguro de deposito"
$ws.Range("G53").Clear()

$ws.Range("A54").Value = 42832
$ws.Range("B54").Value = 500
$ws.Range("C54").Value = 27
$ws.Range("D54").Value = "Medicinas Laila regalo"
$ws.Range("G54").Clear()

$ws.Range("A55").Value = 42832
$ws.Range("B55").Value = -500
$ws.Range("C55").Value = 27
$ws.Range("D55").Value = "Medicinas Laila regalo"
$ws.Range("G55").Clear()

# The frozen-pane view scrolled down to keep the new rows in sight and the last
# selection moved to D40.
$ws.Range("D40").Select()
